$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# Paragraph 3: "{% foreach field in fields.Group %} " -> "{% for field in fields.Group %} "
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
if ($p3.Range.Text -eq "{% foreach field in fields.Group %} `r") {
    $xml3 = '<w:p xmlns:w="' + $wNs + '"><w:r><w:t xml:space="preserve">{% for field in fields.Group %} </w:t></w:r></w:p>'
    $p3.Range.InsertXML($xml3)
}

# ---------------------------------------------------------------------------
# Paragraph 5: paragraph that only holds the "_GoBack" bookmark -> becomes a
# plain empty paragraph (the bookmark is relocated into paragraph 11, below).
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
if ($p5.Range.Text -eq "`r") {
    $xml5 = '<w:p xmlns:w="' + $wNs + '"/>'
    $p5.Range.InsertXML($xml5)
}

# ---------------------------------------------------------------------------
# Paragraph 6: "{% endeach %} " -> "{% end" | "for" | "%}" | " " (4 runs,
# net text becomes "{% endfor%} ", i.e. no space before the closing "%}").
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
if ($p6.Range.Text -eq "{% endeach %} `r") {
    $xml6 = '<w:p xmlns:w="' + $wNs + '">' +
            '<w:r><w:t>{% end</w:t></w:r>' +
            '<w:r><w:t>for</w:t></w:r>' +
            '<w:r><w:t>%}</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
            '</w:p>'
    $p6.Range.InsertXML($xml6)
}

# ---------------------------------------------------------------------------
# Paragraph 10: "{% foreach " | "field " | "in fields.Group %}" | " "
#            -> "{% for" | " " | "field " | "in fields.Group %}" | " "
# (same net text, but "foreach " is now split into "for" + " ").
# ---------------------------------------------------------------------------
$p10 = $d.Paragraphs.Item(10)
if ($p10.Range.Text -eq "{% foreach field in fields.Group %} `r") {
    $xml10 = '<w:p xmlns:w="' + $wNs + '">' +
             '<w:r><w:t>{% for</w:t></w:r>' +
             '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
             '<w:r><w:t xml:space="preserve">field </w:t></w:r>' +
             '<w:r><w:t>in fields.Group %}</w:t></w:r>' +
             '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
             '</w:p>'
    $p10.Range.InsertXML($xml10)
}

# ---------------------------------------------------------------------------
# Paragraph 11: "{{ field.age }}" | " Something else {% endeach %}"
#            -> "{{ field.age }}" | " Something else {% end" | "for " |
#               [bookmark _GoBack] | "%}"
# (net text "{{ field.age }} Something else {% endfor %}"; the _GoBack
# bookmark that used to live alone in paragraph 5 now sits here).
# ---------------------------------------------------------------------------
$p11 = $d.Paragraphs.Item(11)
if ($p11.Range.Text -eq "{{ field.age }} Something else {% endeach %}`r") {
    $xml11 = '<w:p xmlns:w="' + $wNs + '">' +
             '<w:r><w:t>{{ field.age }}</w:t></w:r>' +
             '<w:r><w:t xml:space="preserve"> Something else {% end</w:t></w:r>' +
             '<w:r><w:t xml:space="preserve">for </w:t></w:r>' +
             '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
             '<w:bookmarkEnd w:id="0"/>' +
             '<w:r><w:t>%}</w:t></w:r>' +
             '</w:p>'
    $p11.Range.InsertXML($xml11)
}
